{"js": "// Helper: search for a unique, exact run of text in the document body and\n// replace it with new text in one shot, preserving the formatting of the\n// matched run(s).\nasync function replaceOnce(context, body, findText, newText) {\n  const results = body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1. \"Fall 2014\" -> \"Fall 2015\"\nawait replaceOnce(context, body, \"Fall 2014\", \"Fall 2015\");\n\n// 2. Course description paragraph: rewrite wording (\"Principles...\" /\n//    \"With applications...\" -> \"Covers principles...\" / \"Includes\n//    applications...\") and change \"Compton\" to lower-case \"compton\".\nawait replaceOnce(\n  context,\n  body,\n  \" Principles of energy transfer by radiation. Elements of classical and quantum theory of photon emission; bremsstrahlung, synchrotron radiation. Compton scattering, plasma effects, atomic and molecular electromagnetic transitions. \",\n  \" Covers principles of energy transfer by radiation; elements of classical and quantum theory of photon emission; bremsstrahlung, synchrotron radiation; compton scattering, plasma effects, atomic and molecular electromagnetic transitions. \"\n);\nawait replaceOnce(\n  context,\n  body,\n  \"With applications to current research into astrophysical phenomena.\",\n  \"Includes applications to current research into astrophysical phenomena.\"\n);\n\n// 3. Class time/location: \"10:00 - 11:30, Hearst Field Annex B1\" ->\n//    \"2:00 - 3:30, Campbell Hall 233\"\nawait replaceOnce(\n  context,\n  body,\n  \" 10:00 - 11:30, Hearst Field Annex B1\",\n  \" 2:00 - 3:30, Campbell Hall 233\"\n);\n\n// 4. Office location: \"Hearst Field Annex B54\" -> \"Campbell 455 (or 425 lab)\"\nawait replaceOnce(\n  context,\n  body,\n  \": Hearst Field Annex B54\",\n  \": Campbell 455 (or 425 lab)\"\n);\n\n// 5. Garbled \"Ay250  P9ls4R*@\" -> \"C207  photontrain\"\nawait replaceOnce(\n  context,\n  body,\n  \"  Ay250  P9ls4R*@\",\n  \"  C207  photontrain\"\n);\n\n// 6. \"NO CLASS\" dates and added note about makeup sessions.\nawait replaceOnce(\n  context,\n  body,\n  \"O CLASS 11/11, 11/25, 11/27\",\n  \"O CLASS 9/29, 10/1, 11/19, 11/24.  We will need to schedule 1 additional class and 1 review session.\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $found = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 1. \"Fall 2014\" -> \"Fall 2015\"\nReplace-Text \"Fall 2014\" \"Fall 2015\"\n\n# 2. Course description paragraph: rewrite wording (\"Principles...\" /\n#    \"With applications...\" -> \"Covers principles...\" / \"Includes\n#    applications...\") and change \"Compton\" to lower-case \"compton\".\nReplace-Text \" Principles of energy transfer by radiation. Elements of classical and quantum theory of photon emission; bremsstrahlung, synchrotron radiation. Compton scattering, plasma effects, atomic and molecular electromagnetic transitions. \" \" Covers principles of energy transfer by radiation; elements of classical and quantum theory of photon emission; bremsstrahlung, synchrotron radiation; compton scattering, plasma effects, atomic and molecular electromagnetic transitions. \"\nReplace-Text \"With applications to current research into astrophysical phenomena.\" \"Includes applications to current research into astrophysical phenomena.\"\n\n# 3. Class time/location: \"10:00 - 11:30, Hearst Field Annex B1\" ->\n#    \"2:00 - 3:30, Campbell Hall 233\"\nReplace-Text \" 10:00 - 11:30, Hearst Field Annex B1\" \" 2:00 - 3:30, Campbell Hall 233\"\n\n# 4. Office location: \"Hearst Field Annex B54\" -> \"Campbell 455 (or 425 lab)\"\nReplace-Text \": Hearst Field Annex B54\" \": Campbell 455 (or 425 lab)\"\n\n# 5. Garbled \"Ay250  P9ls4R*@\" -> \"C207  photontrain\"\nReplace-Text \"  Ay250  P9ls4R*@\" \"  C207  photontrain\"\n\n# 6. \"NO CLASS\" dates and added note about makeup sessions.\nReplace-Text \"O CLASS 11/11, 11/25, 11/27\" \"O CLASS 9/29, 10/1, 11/19, 11/24.  We will need to schedule 1 additional class and 1 review session.\"\n"}
